# Refresh the "cryptos" price list (GitHub Actions scheduled update).
# For every changed row we rewrite the Price (D) and/or Volume(1h) (E) cell.
# Where a Price cell is written, its NumberFormat is first forced to "@"
# (Text) so Excel's automatic type-inference doesn't convert numeric-looking
# values (e.g. "647.97", "1.00") into real numbers, which would silently
# drop the source's text formatting (e.g. "1.00" -> 1). Only cells whose
# value actually changes are touched, so every other cell (and its style)
# is left exactly as it was.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-24: price / volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.875.97"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.688.74"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "647.97"
$ws.Range("E5").Value = "  -4.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.88"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.19"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000232"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.313.58"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.73"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.681.42"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.854.78"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.118"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.97"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.37"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.11"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.03"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.837.20"
$ws.Range("E24").Value = "  -0.15%  "

# Rows 25-26: Dai and PEPE swap position (row25 becomes Dai, row26 becomes PEPE)
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000127"
$ws.Range("E26").Value = "  +0.82%  "

# Rows 27-51: remaining price / volume updates
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.12"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.65"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.74"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.686.66"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.43"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "180.03"
$ws.Range("E39").Value = "  +7.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.89"
$ws.Range("E40").Value = "  -5.50%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0903"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.933"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.35"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.57"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000274"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.86"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -3.65%  "
